$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$origStyle = $ws.Range("A1").Style

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "G2" "9"
Set-TextValue "D3" "21.91"
Set-TextValue "G3" "9"
Set-TextValue "D4" "5.386"
Set-TextValue "G4" "9"
Set-TextValue "D5" "0.05808"
Set-TextValue "G5" "9"
Set-TextValue "D6" "3.378"
Set-TextValue "G6" "9"
Set-TextValue "D7" "6.331"
Set-TextValue "G7" "9"
Set-TextValue "D8" "0.8082"
Set-TextValue "G8" "9"
Set-TextValue "D9" "0.9826"
Set-TextValue "G9" "9"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1425"
Set-TextValue "E10" "9WazirXWRX"
Set-TextValue "G10" "9"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07482"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"
Set-TextValue "G11" "9"
Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03208"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "G12" "9"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03054"
Set-TextValue "E13" "12BitrueCoinBTR"
Set-TextValue "G13" "9"
Set-TextValue "B14" "MCDex"
Set-TextValue "C14" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.177"
Set-TextValue "E14" "13MCDexMCB"
Set-TextValue "G14" "9"
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09400"
Set-TextValue "E15" "14BitMartTokenBMX"
Set-TextValue "G15" "9"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001592"
Set-TextValue "E16" "15BitForexTokenBF"
Set-TextValue "G16" "9"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04808"
Set-TextValue "E17" "16CoinExTokenCET"
Set-TextValue "G17" "9"
Set-TextValue "B18" "One"
Set-TextValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005901"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "G18" "9"
Set-TextValue "D19" "0.006239"
Set-TextValue "G19" "9"
Set-TextValue "E20" "19HotbitTokenHTB"
Set-TextValue "G20" "9"
Set-TextValue "D21" "0.0009974"
Set-TextValue "G21" "9"
Set-TextValue "G22" "9"
Set-TextValue "D23" "3.700"
Set-TextValue "G23" "9"
Set-TextValue "D24" "2.245"
Set-TextValue "G24" "9"
Set-TextValue "D25" "0.3203"
Set-TextValue "G25" "9"
Set-TextValue "G26" "9"
Set-TextValue "D27" "0.0003593"
Set-TextValue "E27" "26UpBotsUBXTBestin24h"
Set-TextValue "G27" "9"
Set-TextValue "G28" "9"
Set-TextValue "G29" "9"
Set-TextValue "G30" "9"
Set-TextValue "G31" "9"
Set-TextValue "G32" "9"
Set-TextValue "G33" "9"
Set-TextValue "G34" "9"
Set-TextValue "G35" "9"
Set-TextValue "G36" "9"
Set-TextValue "G37" "9"
Set-TextValue "G38" "9"
Set-TextValue "G39" "9"
Set-TextValue "D40" "0.03885"
Set-TextValue "G40" "9"
Set-TextValue "D41" "0.006371"
Set-TextValue "G41" "9"
Set-TextValue "D42" "0.1072"
Set-TextValue "G42" "9"
Set-TextValue "G43" "9"
Set-TextValue "D44" "0.006686"
Set-TextValue "G44" "9"
Set-TextValue "D45" "0.00005592"
Set-TextValue "G45" "9"
Set-TextValue "G46" "9"
Set-TextValue "D47" "0.3900"
Set-TextValue "G47" "9"
Set-TextValue "D48" "0.1461"
Set-TextValue "G48" "9"
Set-TextValue "G49" "9"
Set-TextValue "G50" "9"
Set-TextValue "G51" "9"
